# A new weekly price observation (2021-12-10) is inserted at row 58 of the
# "Ajo" (Garlic) price table. Excel's native row-insert behavior shifts every
# existing row from 58..163 down by one (to 59..164), carrying over their
# original values/formatting untouched, and the new row 58 is then filled in
# with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 58; this shifts rows 58:163 down to 59:164
# and keeps the sheet's dimension (A1:R163 -> A1:R164) in sync automatically.
$ws.Rows.Item(58).Insert()

# Populate the newly-inserted row 58 with the new weekly record.
$ws.Cells.Item(58, 1).Value2  = 7
$ws.Cells.Item(58, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value2  = "Ñuble"
$ws.Cells.Item(58, 4).Value2  = 44540
$ws.Cells.Item(58, 5).Value2  = 16
$ws.Cells.Item(58, 6).Value2  = 100112003
$ws.Cells.Item(58, 7).Value2  = "Ajo"
$ws.Cells.Item(58, 8).Value2  = "Chino"
$ws.Cells.Item(58, 9).Value2  = "Primera"
$ws.Cells.Item(58, 10).Value2 = 60
$ws.Cells.Item(58, 11).Value2 = 18000
$ws.Cells.Item(58, 12).Value2 = 19000
$ws.Cells.Item(58, 13).Value2 = 18500
$ws.Cells.Item(58, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(58, 15).Value2 = "China"
$ws.Cells.Item(58, 16).Value2 = 1850
$ws.Cells.Item(58, 17).Value2 = 10
$ws.Cells.Item(58, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the existing date number format, matching
# the rest of column D (Excel's row-insert already carries this over from
# row 57 onto the blank new row, but set it explicitly to be safe).
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
